$p = $ppt.ActivePresentation
$tm = $p.TitleMaster
Write-Output ("type=" + $tm.GetType())
$cs = $tm.ColorScheme
for ($i=1; $i -le $cs.Count; $i++) {
    Write-Output ($i.ToString() + " => " + $cs.Colors($i).RGB)
}
